$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The workers' database was updated: the record for EDUARDO CASTRO BLANCO
# (doc 92450475) now comes before the record for PABLO CARABALLO ROMERO
# (doc 73005538), so swap the Doc Trabajador / Nombre Trabajador values
# shown in the two existing data rows.
$ws.Range("C16").Value = "92450475"
$ws.Range("D16").Value = "EDUARDO CASTRO BLANCO"
$ws.Range("C17").Value = "73005538"
$ws.Range("D17").Value = "PABLO CARABALLO ROMERO"
